$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 23.128865
$ws.Cells.Item(2, 8).Value = 69.386595
$ws.Cells.Item(2, 9).Value = 0.7917836846260858
$ws.Cells.Item(2, 10).Value = 0.7917836846260858
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 57.77686466666668
$ws.Cells.Item(2, 14).Value = 173.330594
$ws.Cells.Item(2, 15).Value = 0.5795494107546101
$ws.Cells.Item(2, 16).Value = 0.57954941075461
$ws.Cells.Item(2, 17).Value = 1336.313302998604
$ws.Cells.Item(2, 18).Value = 12026.81972698743
$ws.Cells.Item(2, 19).Value = 0.458877767870162
$ws.Cells.Item(2, 20).Value = 0.458877767870162

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 23.128865
$ws.Cells.Item(3, 8).Value = 69.386595
$ws.Cells.Item(3, 9).Value = 0.7917836846260858
$ws.Cells.Item(3, 10).Value = 0.7917836846260858
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 9.278736333333333
$ws.Cells.Item(3, 14).Value = 27.836209
$ws.Cells.Item(3, 15).Value = 0.09307334701450438
$ws.Cells.Item(3, 16).Value = 0.09307334701450438
$ws.Cells.Item(3, 17).Value = 214.6066400242617
$ws.Cells.Item(3, 18).Value = 1931.459760218355
$ws.Cells.Item(3, 19).Value = 0.07369395763962659
$ws.Cells.Item(3, 20).Value = 0.07369395763962659

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 23.128865
$ws.Cells.Item(4, 8).Value = 69.386595
$ws.Cells.Item(4, 9).Value = 0.7917836846260858
$ws.Cells.Item(4, 10).Value = 0.7917836846260858
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 32.637132
$ws.Cells.Item(4, 14).Value = 97.91139600000001
$ws.Cells.Item(4, 15).Value = 0.3273772422308855
$ws.Cells.Item(4, 16).Value = 0.3273772422308855
$ws.Cells.Item(4, 17).Value = 754.85982001518
$ws.Cells.Item(4, 18).Value = 6793.73838013662
$ws.Cells.Item(4, 19).Value = 0.2592119591162972
$ws.Cells.Item(4, 20).Value = 0.2592119591162972

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.158105666666667
$ws.Cells.Item(5, 8).Value = 3.474317
$ws.Cells.Item(5, 9).Value = 0.03964609469334889
$ws.Cells.Item(5, 10).Value = 0.03964609469334889
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 57.77686466666668
$ws.Cells.Item(5, 14).Value = 173.330594
$ws.Cells.Item(5, 15).Value = 0.5795494107546101
$ws.Cells.Item(5, 16).Value = 0.57954941075461
$ws.Cells.Item(5, 17).Value = 66.91171437269979
$ws.Cells.Item(5, 18).Value = 602.2054293542981
$ws.Cells.Item(5, 19).Value = 0.02297687081825182
$ws.Cells.Item(5, 20).Value = 0.02297687081825182

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1.158105666666667
$ws.Cells.Item(6, 8).Value = 3.474317
$ws.Cells.Item(6, 9).Value = 0.03964609469334889
$ws.Cells.Item(6, 10).Value = 0.03964609469334889
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 9.278736333333333
$ws.Cells.Item(6, 14).Value = 27.836209
$ws.Cells.Item(6, 15).Value = 0.09307334701450438
$ws.Cells.Item(6, 16).Value = 0.09307334701450438
$ws.Cells.Item(6, 17).Value = 10.74575712713922
$ws.Cells.Item(6, 18).Value = 96.711814144253
$ws.Cells.Item(6, 19).Value = 0.003689994729163962
$ws.Cells.Item(6, 20).Value = 0.003689994729163962

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1.158105666666667
$ws.Cells.Item(7, 8).Value = 3.474317
$ws.Cells.Item(7, 9).Value = 0.03964609469334889
$ws.Cells.Item(7, 10).Value = 0.03964609469334889
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 32.637132
$ws.Cells.Item(7, 14).Value = 97.91139600000001
$ws.Cells.Item(7, 15).Value = 0.3273772422308855
$ws.Cells.Item(7, 16).Value = 0.3273772422308855
$ws.Cells.Item(7, 17).Value = 37.797247512948
$ws.Cells.Item(7, 18).Value = 340.175227616532
$ws.Cells.Item(7, 19).Value = 0.01297922914593311
$ws.Cells.Item(7, 20).Value = 0.01297922914593311

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 4.924119999999999
$ws.Cells.Item(8, 8).Value = 14.77236
$ws.Cells.Item(8, 9).Value = 0.1685702206805652
$ws.Cells.Item(8, 10).Value = 0.1685702206805652
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 57.77686466666668
$ws.Cells.Item(8, 14).Value = 173.330594
$ws.Cells.Item(8, 15).Value = 0.5795494107546101
$ws.Cells.Item(8, 16).Value = 0.57954941075461
$ws.Cells.Item(8, 17).Value = 284.5002148424267
$ws.Cells.Item(8, 18).Value = 2560.50193358184
$ws.Cells.Item(8, 19).Value = 0.09769477206619617
$ws.Cells.Item(8, 20).Value = 0.09769477206619614

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 4.924119999999999
$ws.Cells.Item(9, 8).Value = 14.77236
$ws.Cells.Item(9, 9).Value = 0.1685702206805652
$ws.Cells.Item(9, 10).Value = 0.1685702206805652
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 9.278736333333333
$ws.Cells.Item(9, 14).Value = 27.836209
$ws.Cells.Item(9, 15).Value = 0.09307334701450438
$ws.Cells.Item(9, 16).Value = 0.09307334701450438
$ws.Cells.Item(9, 17).Value = 45.68961115369333
$ws.Cells.Item(9, 18).Value = 411.20650038324
$ws.Cells.Item(9, 19).Value = 0.01568939464571383
$ws.Cells.Item(9, 20).Value = 0.01568939464571383

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 4.924119999999999
$ws.Cells.Item(10, 8).Value = 14.77236
$ws.Cells.Item(10, 9).Value = 0.1685702206805652
$ws.Cells.Item(10, 10).Value = 0.1685702206805652
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 32.637132
$ws.Cells.Item(10, 14).Value = 97.91139600000001
$ws.Cells.Item(10, 15).Value = 0.3273772422308855
$ws.Cells.Item(10, 16).Value = 0.3273772422308855
$ws.Cells.Item(10, 17).Value = 160.70915442384
$ws.Cells.Item(10, 18).Value = 1446.38238981456
$ws.Cells.Item(10, 19).Value = 0.05518605396865524
$ws.Cells.Item(10, 20).Value = 0.05518605396865524
